$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value would otherwise be auto-parsed as a number by Excel;
# force Text format first so the literal string is preserved (matches original inline-string cells),
# then clear the format footprint so no stray style is left behind.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '62.215.20'
$ws.Range('E2').Value = '  -2.82%  '
$ws.Range('D3').Value = '2.985.95'
$ws.Range('E3').Value = '  -3.53%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '583.64'
$ws.Range('E5').Value = '  -2.27%  '
$ws.Range('D6').Value = '145.24'
$ws.Range('E6').Value = '  -7.02%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -3.52%  '
$ws.Range('D9').Value = '2.985.75'
$ws.Range('E10').Value = '  -7.09%  '
$ws.Range('D11').Value = '5.71'
$ws.Range('E11').Value = '  -4.36%  '
$ws.Range('E12').Value = '  -2.73%  '
$ws.Range('E13').Value = '  -5.38%  '
$ws.Range('D14').Value = '34.40'
$ws.Range('E14').Value = '  -6.79%  '
$ws.Range('E15').Value = '  +1.78%  '
$ws.Range('D16').Value = '3.476.95'
$ws.Range('E16').Value = '  -3.59%  '
$ws.Range('D17').Value = '62.177.26'
$ws.Range('E17').Value = '  -2.72%  '
$ws.Range('D18').Value = '6.95'
$ws.Range('E18').Value = '  -3.73%  '
$ws.Range('D19').Value = '2.988.85'
$ws.Range('E19').Value = '  -3.41%  '
$ws.Range('D20').Value = '455.44'
$ws.Range('E20').Value = '  -5.30%  '
$ws.Range('E21').Value = '  -4.85%  '
$ws.Range('E22').Value = '  -5.52%  '
$ws.Range('E23').Value = '  -3.45%  '
$ws.Range('D24').Value = '79.86'
$ws.Range('E24').Value = '  -2.19%  '
$ws.Range('D25').Value = '2.23'
$ws.Range('E25').Value = '  -10.55%  '
$ws.Range('D26').Value = '12.17'
$ws.Range('E26').Value = '  -5.64%  '
$ws.Range('D27').Value = '10.03'
$ws.Range('E27').Value = '  -6.79%  '
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').Value = '2.60'
$ws.Range('E30').Value = '  -3.70%  '
$ws.Range('D31').Value = '7.11'
$ws.Range('E31').Value = '  -7.05%  '
$ws.Range('E32').Value = '  -6.29%  '
$ws.Range('D33').Value = '26.67'
$ws.Range('E33').Value = '  -2.01%  '
$ws.Range('E34').Value = '  -5.01%  '
$ws.Range('E35').Value = '  -5.49%  '
$ws.Range('D36').Value = '0.0₃0780'
$ws.Range('E36').Value = '  -7.48%  '
$ws.Range('E37').Value = '  -5.70%  '
$ws.Range('E38').Value = '  -7.07%  '
$ws.Range('D39').Value = '50.01'
$ws.Range('E39').Value = '  -1.99%  '
$ws.Range('D40').Value = '8.90'
$ws.Range('E40').Value = '  -3.72%  '
$ws.Range('E41').Value = '  -11.68%  '
$ws.Range('D42').Value = '396.23'
$ws.Range('E42').Value = '  -10.84%  '
$ws.Range('E43').Value = '  +0.65%  '
$ws.Range('D44').Value = '0.271'
$ws.Range('E44').Value = '  -7.23%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.748.31'
$ws.Range('E45').Value = '  -3.15%  '
$ws.Range('B46').Value = 'Arweave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D46').Value = '38.75'
$ws.Range('E46').Value = '  -3.38%  '
$ws.Range('E47').Value = '  -4.47%  '
$ws.Range('D48').Value = '127.37'
$ws.Range('E48').Value = '  -3.66%  '
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('E50').Value = '  -2.39%  '
$ws.Range('D51').Value = '23.54'
$ws.Range('E51').Value = '  -9.78%  '

# Remove the temporary Text-format styling footprint from the forced cells above,
# restoring them to the unstyled state while keeping their text cell-type.
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D51").ClearFormats()
